$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 751
$ws1.Range("F8").Value = 6229
$ws1.Range("F12").Value = 5179
$ws1.Range("F15").Value = 1179
$ws1.Range("F22").Value = 3697

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 81

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 81
$ws4.Range("F4").Value = 751
$ws4.Range("F9").Value = 6229
$ws4.Range("F13").Value = 5179
$ws4.Range("F16").Value = 1179
$ws4.Range("F23").Value = 3697
